$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the formatting of the other
# header cells (e.g. G1: bold, centered, bordered) by copying that cell's
# format rather than re-creating it (keeps the existing style index).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new data value in H2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 0
